# تعديل تلقائي في شيت Card21 by admin at 2025-12-10 13:06:55
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# D8: "nan" -> "966.5" (must stay a literal text value, matching the rest
# of the sheet's numeric-looking-but-textual cells). Format as Text first so
# Excel doesn't auto-convert the numeric-looking string to a Number, then
# clear the formatting footprint so the cell keeps its original (default)
# style while the stored value remains text.
$d8 = $ws.Cells.Item(8, 4)
$d8.NumberFormat = "@"
$d8.Value = "966.5"
$d8.ClearFormats()

# Row 29, columns B:K were emptied inlineStr cells; restore them to the
# literal text "nan" (matching sibling rows such as row 8/row 2 in the same
# columns).
$row29Cols = 2..11   # B=2 .. K=11
foreach ($col in $row29Cols) {
    $ws.Cells.Item(29, $col).Value = "nan"
}
